$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Coin/Link/Price/Volume columns so numeric-
# looking strings (e.g. "1.000", "22.90") keep their literal text
# representation instead of being coerced to a Number by COM.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.003.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.88%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.828.40'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.83%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.16%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6528'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.20%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.49'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.70%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2930'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.72%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07327'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.90'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.11%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07665'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.68%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.833.49'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.19%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.974'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.63%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6647'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.13%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.95'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.90%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.055'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.30%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008641'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.76%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '28.906.70'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.27%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.083.41'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.77%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.40'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.93%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.61'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.086'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.54%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.01%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.82'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.87%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.497'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.48%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1377'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.99%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.89'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.506'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.23%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.099'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.86%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.009'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.57%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05341'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.61%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7420'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.829'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.35%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.151'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.31%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.20%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.298.20'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.52%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01783'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.28%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.85%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.360'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +6.66%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8952'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.21%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.86'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.68%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.982.21'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5139'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.50%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.90'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.57%  '

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000120'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.03%  '

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.732'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07301'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -11.37%  '
